# Updates the cryptocurrency price / 1h-volume snapshot (and restores the
# Polygon / WrappedEther row order, which the upstream feed swapped) to match
# the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = 'D2'; Value = '29.926.67' },
    @{ Ref = 'E2'; Value = '  +0.19%  ' },
    @{ Ref = 'D3'; Value = '1.876.20' },
    @{ Ref = 'E3'; Value = '  -0.60%  ' },
    @{ Ref = 'D4'; Value = '1.001' },
    @{ Ref = 'E4'; Value = '  -0.03%  ' },
    @{ Ref = 'D5'; Value = '0.7405' },
    @{ Ref = 'E5'; Value = '  -3.63%  ' },
    @{ Ref = 'D6'; Value = '242.71' },
    @{ Ref = 'E6'; Value = '  +0.04%  ' },
    @{ Ref = 'D7'; Value = '1.002' },
    @{ Ref = 'E7'; Value = '  +0.09%  ' },
    @{ Ref = 'D8'; Value = '0.3156' },
    @{ Ref = 'E8'; Value = '  +1.04%  ' },
    @{ Ref = 'D9'; Value = '0.07229' },
    @{ Ref = 'E9'; Value = '  +0.87%  ' },
    @{ Ref = 'D10'; Value = '24.64' },
    @{ Ref = 'E10'; Value = '  -3.86%  ' },
    @{ Ref = 'D11'; Value = '0.08348' },
    @{ Ref = 'E11'; Value = '  -2.55%  ' },
    @{ Ref = 'B12'; Value = 'WrappedEther' },
    @{ Ref = 'C12'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Ref = 'D12'; Value = '1.911.09' },
    @{ Ref = 'E12'; Value = '  -0.15%  ' },
    @{ Ref = 'B13'; Value = 'Polygon' },
    @{ Ref = 'C13'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' },
    @{ Ref = 'D13'; Value = '0.7503' },
    @{ Ref = 'E13'; Value = '  -1.71%  ' },
    @{ Ref = 'D14'; Value = '5.399' },
    @{ Ref = 'E14'; Value = '  +0.68%  ' },
    @{ Ref = 'D15'; Value = '92.38' },
    @{ Ref = 'E15'; Value = '  -1.31%  ' },
    @{ Ref = 'D16'; Value = '29.948.10' },
    @{ Ref = 'E16'; Value = '  +0.29%  ' },
    @{ Ref = 'D17'; Value = '6.106' },
    @{ Ref = 'E17'; Value = '  -0.77%  ' },
    @{ Ref = 'D18'; Value = '248.19' },
    @{ Ref = 'E18'; Value = '  +1.50%  ' },
    @{ Ref = 'D19'; Value = '13.58' },
    @{ Ref = 'E19'; Value = '  -1.41%  ' },
    @{ Ref = 'D20'; Value = '0.000007844' },
    @{ Ref = 'E20'; Value = '  +0.49%  ' },
    @{ Ref = 'E21'; Value = '  +0.26%  ' },
    @{ Ref = 'D22'; Value = '2.143.75' },
    @{ Ref = 'E22'; Value = '  +0.02%  ' },
    @{ Ref = 'D23'; Value = '8.022' },
    @{ Ref = 'E23'; Value = '  +0.06%  ' },
    @{ Ref = 'E24'; Value = '  -0.07%  ' },
    @{ Ref = 'D25'; Value = '0.1554' },
    @{ Ref = 'E25'; Value = '  -5.16%  ' },
    @{ Ref = 'D26'; Value = '9.279' },
    @{ Ref = 'E26'; Value = '  -1.13%  ' },
    @{ Ref = 'D27'; Value = '164.77' },
    @{ Ref = 'E27'; Value = '  +1.06%  ' },
    @{ Ref = 'D28'; Value = '18.68' },
    @{ Ref = 'E28'; Value = '  -0.27%  ' },
    @{ Ref = 'D29'; Value = '2.028' },
    @{ Ref = 'E29'; Value = '  -0.17%  ' },
    @{ Ref = 'D30'; Value = '1.508' },
    @{ Ref = 'E30'; Value = '  +2.67%  ' },
    @{ Ref = 'D31'; Value = '4.603' },
    @{ Ref = 'E31'; Value = '  +2.12%  ' },
    @{ Ref = 'D32'; Value = '1.538' },
    @{ Ref = 'E32'; Value = '  -0.08%  ' },
    @{ Ref = 'D33'; Value = '4.265' },
    @{ Ref = 'E33'; Value = '  +4.13%  ' },
    @{ Ref = 'D34'; Value = '0.05328' },
    @{ Ref = 'E34'; Value = '  -2.19%  ' },
    @{ Ref = 'D35'; Value = '1.235' },
    @{ Ref = 'E35'; Value = '  -0.42%  ' },
    @{ Ref = 'D36'; Value = '0.7502' },
    @{ Ref = 'E36'; Value = '  +1.03%  ' },
    @{ Ref = 'D37'; Value = '1.002' },
    @{ Ref = 'E37'; Value = '  +0.16%  ' },
    @{ Ref = 'D38'; Value = '2.699' },
    @{ Ref = 'E38'; Value = '  +0.18%  ' },
    @{ Ref = 'E39'; Value = '  +0.65%  ' },
    @{ Ref = 'D40'; Value = '2.754' },
    @{ Ref = 'E40'; Value = '  -1.10%  ' },
    @{ Ref = 'D41'; Value = '0.4547' },
    @{ Ref = 'E41'; Value = '  +1.75%  ' },
    @{ Ref = 'D42'; Value = '6.151' },
    @{ Ref = 'E42'; Value = '  +1.33%  ' },
    @{ Ref = 'D43'; Value = '1.103.41' },
    @{ Ref = 'E43'; Value = '  -0.20%  ' },
    @{ Ref = 'D44'; Value = '72.38' },
    @{ Ref = 'E44'; Value = '  -1.00%  ' },
    @{ Ref = 'D45'; Value = '0.8629' },
    @{ Ref = 'E45'; Value = '  +1.30%  ' },
    @{ Ref = 'D46'; Value = '104.36' },
    @{ Ref = 'E46'; Value = '  +1.69%  ' },
    @{ Ref = 'D47'; Value = '1.002' },
    @{ Ref = 'E47'; Value = '  +0.15%  ' },
    @{ Ref = 'D48'; Value = '1.865' },
    @{ Ref = 'E48'; Value = '  +0.19%  ' },
    @{ Ref = 'D49'; Value = '7.616' },
    @{ Ref = 'E49'; Value = '  -0.60%  ' },
    @{ Ref = 'D50'; Value = '9.526' },
    @{ Ref = 'E50'; Value = '  -2.52%  ' },
    @{ Ref = 'D51'; Value = '2.040.28' },
    @{ Ref = 'E51'; Value = '  -0.21%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $text = $u.Value

    # Columns B/C (coin name / link) and column E (padded "  +x.xx%  " strings)
    # are never ambiguous with numbers, so a plain assignment is safe and keeps
    # the cell's existing (default) style untouched.
    #
    # Column D quotes are plain decimals (e.g. "1.001", "0.7405") that Excel's
    # COM layer would otherwise auto-coerce to a Number on assignment -- losing
    # the literal text the source feed provides (and collapsing values like
    # "1.001"/"1.002" together via float rounding). Force the cell to Text via
    # NumberFormat, assign, then reset the style back to the sheet default so
    # no stray formatting is introduced (matching cells keep their original,
    # unstyled look).
    if ($u.Ref.StartsWith("D")) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}
